$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    5799.800999999999,
    5649.136500000001,
    5554.8325,
    5481.98,
    5456.8815,
    5512.395,
    5765.477227722772,
    6519.4465,
    7901.5895,
    13052.0005,
    14886.7985,
    14810.005,
    14699.3385,
    15357.636,
    15432.207,
    15338.9075,
    15792.80217625723,
    15557.50171551809,
    15152.49877462994,
    13812.90902852661,
    12007.84432898735,
    9523.143,
    6375.7855,
    5494.996500000001
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
